# Updated symbol list on Sun Dec 18 15:07:57 UTC 2022 with GitHub Actions
#
# Refresh the scraped coinranking.com price table: every row's "Hora"
# (hour, col G) ticks from 14 -> 15, many "Price" values (col D) move to
# their newly-scraped quotes, and rows 42/43 (CEJI / BKEXToken) swap
# positions with refreshed prices for their new slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text (Price col D, Hora col G) must stay
# plain text (matching the sheet's existing t="inlineStr" cells), so they
# are written through .Formula with a leading apostrophe to force text
# instead of letting Excel auto-convert them to numbers.
function Set-TextValue($addr, $value) {
    $ws.Range($addr).Formula = "'" + $value
}

# row -> new Price (col D); omitted rows keep their existing price
$priceUpdates = @{
    2  = '246.45'
    4  = '5.466'
    5  = '0.05635'
    6  = '6.468'
    7  = '0.8058'
    8  = '1.047'
    9  = '0.1438'
    10 = '0.07354'
    11 = '0.03202'
    12 = '0.02936'
    13 = '0.09265'
    14 = '0.001676'
    15 = '3.205'
    16 = '0.04735'
    17 = '0.0005850'
    18 = '0.006383'
    19 = '0.001058'
    20 = '0.004113'
    22 = '3.978'
    23 = '3.388'
    24 = '2.129'
    27 = '0.0003010'
    40 = '0.04152'
    41 = '0.006883'
    44 = '0.009029'
    45 = '0.00005661'
    47 = '0.6822'
    48 = '0.01909'
    49 = '0.00002107'
}

# Rows 42 and 43 swap their Coin/Link/Volume identity (CEJI <-> BKEXToken)
# and each gets a freshly scraped price for its new row slot.
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$priceUpdates[42] = '0.1039'

$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E43').Value = '42CEJICEJI'
$priceUpdates[43] = '0.002980'

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue "D$row" $priceUpdates[$row]
}

# Every data row (2-51) advances "Hora" from 14 to 15.
for ($row = 2; $row -le 51; $row++) {
    Set-TextValue "G$row" '15'
}
